# Add a new "maximum" column (F) to the animalBreakpointslistRaw sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("animalBreakpointslistRaw")

# Header
$ws.Range("F1").Value = "maximum"

# Data values for column F (rows 2-8)
$values = @(115, 50, 115, 100, 70, 50, 50)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $values[$i]
    # Match the numeric formatting used by the other measurement columns (e.g. column E)
    $cell.NumberFormat = $ws.Cells.Item($row, 5).NumberFormat
}

# Update the selected cell to reflect the newly entered value, like Excel would
# after typing into F2.
$ws.Range("F2").Select()
